# Revert "feat: update config" — restore the original "bgm" (column U)
# values on Sheet1 that the prior commit had bumped to 434716 across the
# board. Column U = "bgm", grouped by stage id ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("U5:U20").Value = 162450
$ws.Range("U21:U50").Value = 404268
$ws.Range("U51:U80").Value = 404267
$ws.Range("U81:U110").Value = 404269
$ws.Range("U111:U140").Value = 404266
$ws.Range("U141:U170").Value = 404264
$ws.Range("U171").Value = 404270
$ws.Range("U172").Value = 404268

# Rows 21, 50, 80, 110, 140 carry an explicit row-level custom format
# (fill style index 1 / green, same as the other edited cells in those
# rows). Touching their U cell causes Excel to bake that effective
# style onto the cell explicitly too, so replicate it here.
foreach ($r in @(21, 50, 80, 110, 140)) {
    $cell = $ws.Cells.Item($r, 21)
    $cell.Interior.Color = 3385600
    $cell.Interior.PatternColor = 32768
}
